# "Measurements for Components.xlsx" update
# - LED Panel's measurement changes from 81x37 -> 71x23.5
# - Two new rows of data are recorded for an "Enclosure" component
#   (one ends up with a #VALUE! error from the existing helper column
#   formula, the other keeps a broken #REF! formula - exactly mirroring
#   what the author's worksheet looked like after they typed the new
#   rows in over the old helper rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Enclosure" entry typed into what used to be blank helper rows.
$ws.Range("B15").Value = "Enclosure"

# LED Panel dimensions updated
$ws.Range("C9").Value = "71x23.5"

$ws.Range("B16").Value = "251.7x223.8x50.4"

# Row 16's running-count helper formula ends up referencing a cell that
# no longer resolves, leaving a #REF! error baked into the formula text.
$ws.Range("A16").Formula = "=IF(#REF!<>"""", A15+1, """")"

# Leave the cursor where the author finished typing.
$ws.Activate()
[void]$ws.Range("A9").Select()
$excel.ActiveWindow.ScrollRow = 9
[void]$ws.Range("B16").Select()
